$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 39306.31640625
$ws.Range("D3").Value = 39306.31640625

$ws.Range("A4").Value = 2.0
$ws.Range("C4").Value = 2233.0
$ws.Range("D4").Value = 10.0
$ws.Range("A5").Value = 29.0
$ws.Range("C5").Value = 3177.0
$ws.Range("D5").Value = 756.0
$ws.Range("A6").Value = 41.0
$ws.Range("C6").Value = 4985.0
$ws.Range("D6").Value = 140.0
$ws.Range("A7").Value = 34.0
$ws.Range("C7").Value = 4608.0
$ws.Range("D7").Value = 1198.0
$ws.Range("A8").Value = 16.0
$ws.Range("C8").Value = 6107.0
$ws.Range("D8").Value = 669.0
$ws.Range("A9").Value = 22.0
$ws.Range("C9").Value = 6101.0
$ws.Range("D9").Value = 1110.0
$ws.Range("A10").Value = 1.0
$ws.Range("C10").Value = 6734.0
$ws.Range("D10").Value = 1453.0
$ws.Range("A11").Value = 8.0
$ws.Range("C11").Value = 7265.0
$ws.Range("D11").Value = 1268.0
$ws.Range("A12").Value = 38.0
$ws.Range("C12").Value = 7392.0
$ws.Range("D12").Value = 2244.0
$ws.Range("A13").Value = 31.0
$ws.Range("C13").Value = 7545.0
$ws.Range("D13").Value = 2801.0
$ws.Range("A14").Value = 9.0
$ws.Range("C14").Value = 6898.0
$ws.Range("D14").Value = 1885.0
$ws.Range("A15").Value = 40.0
$ws.Range("C15").Value = 6271.0
$ws.Range("D15").Value = 2135.0
$ws.Range("A16").Value = 3.0
$ws.Range("C16").Value = 5530.0
$ws.Range("D16").Value = 1424.0
$ws.Range("A17").Value = 14.0
$ws.Range("C17").Value = 4612.0
$ws.Range("D17").Value = 2035.0
$ws.Range("A18").Value = 23.0
$ws.Range("C18").Value = 5199.0
$ws.Range("D18").Value = 2182.0
$ws.Range("A19").Value = 11.0
$ws.Range("C19").Value = 5468.0
$ws.Range("D19").Value = 2606.0
$ws.Range("A20").Value = 15.0
$ws.Range("C20").Value = 6347.0
$ws.Range("D20").Value = 2683.0
$ws.Range("A21").Value = 12.0
$ws.Range("C21").Value = 5989.0
$ws.Range("D21").Value = 2873.0
$ws.Range("A22").Value = 46.0
$ws.Range("C22").Value = 6807.0
$ws.Range("D22").Value = 2993.0
$ws.Range("A23").Value = 44.0
$ws.Range("C23").Value = 7509.0
$ws.Range("D23").Value = 3239.0
$ws.Range("A24").Value = 18.0
$ws.Range("C24").Value = 7462.0
$ws.Range("D24").Value = 3590.0
$ws.Range("A25").Value = 7.0
$ws.Range("C25").Value = 7573.0
$ws.Range("D25").Value = 3716.0
$ws.Range("A26").Value = 28.0
$ws.Range("C26").Value = 7541.0
$ws.Range("D26").Value = 3981.0
$ws.Range("A27").Value = 6.0
$ws.Range("C27").Value = 7608.0
$ws.Range("D27").Value = 4458.0
$ws.Range("A28").Value = 30.0
$ws.Range("C28").Value = 7352.0
$ws.Range("D28").Value = 4506.0
$ws.Range("A29").Value = 43.0
$ws.Range("C29").Value = 7280.0
$ws.Range("D29").Value = 4899.0
$ws.Range("A30").Value = 17.0
$ws.Range("C30").Value = 7611.0
$ws.Range("D30").Value = 5184.0
$ws.Range("A31").Value = 27.0
$ws.Range("C31").Value = 7555.0
$ws.Range("D31").Value = 4819.0
$ws.Range("A32").Value = 19.0
$ws.Range("C32").Value = 7732.0
$ws.Range("D32").Value = 4723.0
$ws.Range("A33").Value = 37.0
$ws.Range("C33").Value = 7762.0
$ws.Range("D33").Value = 4595.0
$ws.Range("A34").Value = 36.0
$ws.Range("C34").Value = 7248.0
$ws.Range("D34").Value = 3779.0
$ws.Range("A35").Value = 33.0
$ws.Range("C35").Value = 6426.0
$ws.Range("D35").Value = 3173.0
$ws.Range("A36").Value = 20.0
$ws.Range("C36").Value = 5900.0
$ws.Range("D36").Value = 3561.0
$ws.Range("A37").Value = 47.0
$ws.Range("C37").Value = 5185.0
$ws.Range("D37").Value = 3258.0
$ws.Range("A38").Value = 13.0
$ws.Range("C38").Value = 4706.0
$ws.Range("D38").Value = 2674.0
$ws.Range("A39").Value = 21.0
$ws.Range("C39").Value = 4483.0
$ws.Range("D39").Value = 3369.0
$ws.Range("A40").Value = 25.0
$ws.Range("C40").Value = 4307.0
$ws.Range("D40").Value = 2322.0
$ws.Range("A41").Value = 5.0
$ws.Range("C41").Value = 3082.0
$ws.Range("D41").Value = 1644.0
$ws.Range("A42").Value = 48.0
$ws.Range("C42").Value = 3023.0
$ws.Range("D42").Value = 1942.0
$ws.Range("A43").Value = 39.0
$ws.Range("C43").Value = 3484.0
$ws.Range("D43").Value = 2829.0
$ws.Range("A44").Value = 32.0
$ws.Range("C44").Value = 3245.0
$ws.Range("D44").Value = 3305.0
$ws.Range("A45").Value = 42.0
$ws.Range("C45").Value = 1916.0
$ws.Range("D45").Value = 1569.0
$ws.Range("A46").Value = 24.0
$ws.Range("C46").Value = 1633.0
$ws.Range("D46").Value = 2809.0
$ws.Range("A47").Value = 10.0
$ws.Range("C47").Value = 1112.0
$ws.Range("D47").Value = 2049.0
$ws.Range("A48").Value = 45.0
$ws.Range("C48").Value = 10.0
$ws.Range("D48").Value = 2676.0
$ws.Range("A49").Value = 35.0
$ws.Range("C49").Value = 23.0
$ws.Range("D49").Value = 2216.0
$ws.Range("A50").Value = 26.0
$ws.Range("C50").Value = 675.0
$ws.Range("D50").Value = 1006.0
$ws.Range("A51").Value = 4.0
$ws.Range("C51").Value = 401.0
$ws.Range("D51").Value = 841.0
